$d = $word.ActiveDocument

# The second paragraph in the document (currently empty) needs to be
# turned into five paragraphs containing "Fhk", "1", "1", "2", "3".
$target = $d.Paragraphs.Item(2)

# Insert the new text (with paragraph breaks) right before the existing
# paragraph mark of that empty paragraph. This turns the single empty
# paragraph into: Fhk / 1 / 1 / 2 / 3 / (original empty paragraph mark).
$target.Range.InsertBefore("Fhk`r1`r1`r2`r3")
